$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("M2").Value = 22.618885
$ws.Range("N2").Value = 67.856655
$ws.Range("O2").Value = 0.9027998993061069
$ws.Range("P2").Value = 0.902799899306107
$ws.Range("Q2").Value = 3.616767251128334
$ws.Range("R2").Value = 32.550905260155
$ws.Range("S2").Value = 0.02390408243948202
$ws.Range("T2").Value = 0.02390408243948203

# Row 3
$ws.Range("O3").Value = 0.08600330007856447
$ws.Range("P3").Value = 0.08600330007856449
$ws.Range("S3").Value = 0.002277171249936591
$ws.Range("T3").Value = 0.002277171249936592

# Row 4
$ws.Range("M4").Value = 0.2805263333333333
$ws.Range("N4").Value = 0.841579
$ws.Range("O4").Value = 0.0111968006153285
$ws.Range("P4").Value = 0.01119680061532851
$ws.Range("Q4").Value = 0.04485625420877778
$ws.Range("R4").Value = 0.403706287879
$ws.Range("S4").Value = 0.0002964657452586904
$ws.Range("T4").Value = 0.0002964657452586905

# Row 5
$ws.Range("M5").Value = 22.618885
$ws.Range("N5").Value = 67.856655
$ws.Range("O5").Value = 0.9027998993061069
$ws.Range("P5").Value = 0.902799899306107
$ws.Range("Q5").Value = 132.97986298552
$ws.Range("R5").Value = 1196.81876686968
$ws.Range("S5").Value = 0.8788958168666249
$ws.Range("T5").Value = 0.8788958168666251

# Row 6
$ws.Range("O6").Value = 0.08600330007856447
$ws.Range("P6").Value = 0.08600330007856449
$ws.Range("S6").Value = 0.08372612882862787
$ws.Range("T6").Value = 0.08372612882862791

# Row 7
$ws.Range("M7").Value = 0.2805263333333333
$ws.Range("N7").Value = 0.841579
$ws.Range("O7").Value = 0.0111968006153285
$ws.Range("P7").Value = 0.01119680061532851
$ws.Range("Q7").Value = 1.649256953669333
$ws.Range("R7").Value = 14.843312583024
$ws.Range("S7").Value = 0.01090033487006981
$ws.Range("T7").Value = 0.01090033487006982
